$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume snapshot (prices, % changes) and
# re-rank three coins (TrustWalletToken now above TheSandbox above Quant).
# Price cells that look like plain numbers are forced to Text format first
# so Excel keeps them as literal strings (e.g. "1.001") instead of coercing
# them into floating-point numbers.

$ws.Range('D2').Value = '30.977.75'
$ws.Range('E2').Value = '  +1.11%  '
$ws.Range('D3').Value = '1.953.83'
$ws.Range('E3').Value = '  -0.30%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.33'
$ws.Range('E5').Value = '  -1.34%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4898'
$ws.Range('E7').Value = '  +1.34%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2953'
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06822'
$ws.Range('E9').Value = '  +0.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.15'
$ws.Range('E10').Value = '  -0.87%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '107.04'
$ws.Range('E11').Value = '  -3.29%  '
$ws.Range('D12').Value = '1.952.93'
$ws.Range('E12').Value = '  -0.43%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07779'
$ws.Range('E13').Value = '  +0.67%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.421'
$ws.Range('E14').Value = '  -1.10%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.7017'
$ws.Range('E15').Value = '  +1.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '281.91'
$ws.Range('E16').Value = '  -3.70%  '
$ws.Range('D17').Value = '31.009.55'
$ws.Range('E17').Value = '  +1.16%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.18'
$ws.Range('E18').Value = '  -1.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007670'
$ws.Range('E19').Value = '  -0.33%  '
$ws.Range('D20').Value = '2.210.30'
$ws.Range('E20').Value = '  -0.18%  '
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.475'
$ws.Range('E22').Value = '  -3.31%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.474'
$ws.Range('E24').Value = '  -2.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.809'
$ws.Range('E25').Value = '  -0.54%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '169.69'
$ws.Range('E26').Value = '  -0.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.97'
$ws.Range('E27').Value = '  -0.93%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.195'
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1054'
$ws.Range('E29').Value = '  -1.88%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.420'
$ws.Range('E30').Value = '  -1.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.579'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.615'
$ws.Range('E32').Value = '  -0.89%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.424'
$ws.Range('E33').Value = '  -0.38%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04930'
$ws.Range('E34').Value = '  -3.26%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7634'
$ws.Range('E35').Value = '  -1.94%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.167'
$ws.Range('E36').Value = '  -0.67%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.727'
$ws.Range('E37').Value = '  -0.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02008'
$ws.Range('E38').Value = '  -2.49%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.703'
$ws.Range('E39').Value = '  -0.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.540'
$ws.Range('E40').Value = '  +7.35%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.124'
$ws.Range('E41').Value = '  +2.62%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '73.81'
$ws.Range('E42').Value = '  +5.02%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8844'
$ws.Range('E43').Value = '  +0.98%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4460'
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '109.24'
$ws.Range('E45').Value = '  -1.71%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.085'
$ws.Range('E46').Value = '  +9.26%  '
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '992.93'
$ws.Range('E48').Value = '  +8.99%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1259'
$ws.Range('E49').Value = '  -2.33%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.320'
$ws.Range('E50').Value = '  -0.26%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.2571'
$ws.Range('E51').Value = '  +2.19%  '
